$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 645, shifting existing rows 645:743 down to 646:744.
$ws.Rows.Item(645).Insert()

# Populate the newly inserted row 645 with the new price-report entry.
$ws.Range("A645").Value = 7
$ws.Range("B645").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C645").Value = "Ñuble"
$ws.Range("D645").Value = 45218
$ws.Range("E645").Value = 16
$ws.Range("F645").Value = 100114001
$ws.Range("G645").Value = "Papa"
$ws.Range("H645").Value = "Asterix"
$ws.Range("I645").Value = "1a (guarda)"
$ws.Range("J645").Value = 200
$ws.Range("K645").Value = 28000
$ws.Range("L645").Value = 28000
$ws.Range("M645").Value = 28000
$ws.Range("N645").Value = "$/saco 25 kilos"
$ws.Range("O645").Value = "Región de Los Lagos"
$ws.Range("P645").Value = 1120
$ws.Range("Q645").Value = 25
$ws.Range("R645").Value = "Hortaliza"
